# PRIXMM.xlsx — "Add files via upload" re-upload.
# Observed changes: four quantity/price figures in column D were revised
# upward, and the sheet's scroll/selection state moved from the bottom of
# the list (row 36 / C36) up to around row 19-29 with D29 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Revised figures in column D
$ws.Range("D17").Value = 331   # was 302
$ws.Range("D20").Value = 284   # was 278
$ws.Range("D22").Value = 168   # was 164.5
$ws.Range("D29").Value = 130   # was 112.5

# Restore the scrolled/selected view state (top-left around row 19,
# D29 as the active cell) as closely as the object model allows.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D29").Select()
